$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$lo = $ws.ListObjects.Item(1)

# Extend Table1 ("Harian 23.001") with new biweekly "upah harian" rows (Jul-Nov 2023)
# by adding table rows, copying the format of the row above, then writing the new data.
$lo.ListRows.Add() | Out-Null
$ws.Range("B41:E41").Copy()
$ws.Range("B42:E42").PasteSpecial(-4122)
$ws.Range("B42").Value = 45116
$ws.Range("C42").Value = 45129
$ws.Range("D42").Value = "Harian Driver"
$ws.Range("E42").Value = 24225400

$lo.ListRows.Add() | Out-Null
$ws.Range("B42:E42").Copy()
$ws.Range("B43:E43").PasteSpecial(-4122)
$ws.Range("B43").Value = 45130
$ws.Range("C43").Value = 45143
$ws.Range("D43").Value = "Harian Driver"
$ws.Range("E43").Value = 22271200

$lo.ListRows.Add() | Out-Null
$ws.Range("B43:E43").Copy()
$ws.Range("B44:E44").PasteSpecial(-4122)
$ws.Range("B44").Value = 45144
$ws.Range("C44").Value = 45157
$ws.Range("D44").Value = "Harian Driver"
$ws.Range("E44").Value = 19299400

$lo.ListRows.Add() | Out-Null
$ws.Range("B44:E44").Copy()
$ws.Range("B45:E45").PasteSpecial(-4122)
$ws.Range("B45").Value = 45158
$ws.Range("C45").Value = 45171
$ws.Range("D45").Value = "Harian Driver"
$ws.Range("E45").Value = 21150000

$lo.ListRows.Add() | Out-Null
$ws.Range("B45:E45").Copy()
$ws.Range("B46:E46").PasteSpecial(-4122)
$ws.Range("B46").Value = 45172
$ws.Range("C46").Value = 45185
$ws.Range("D46").Value = "Harian Driver"
$ws.Range("E46").Value = 18925200

$lo.ListRows.Add() | Out-Null
$ws.Range("B46:E46").Copy()
$ws.Range("B47:E47").PasteSpecial(-4122)
$ws.Range("B47").Value = 45186
$ws.Range("C47").Value = 45199
$ws.Range("D47").Value = "Harian Driver"
$ws.Range("E47").Value = 19388400

$lo.ListRows.Add() | Out-Null
$ws.Range("B47:E47").Copy()
$ws.Range("B48:E48").PasteSpecial(-4122)
$ws.Range("B48").Value = 45200
$ws.Range("C48").Value = 45213
$ws.Range("D48").Value = "Harian Driver"
$ws.Range("E48").Value = 21710200

$lo.ListRows.Add() | Out-Null
$ws.Range("B48:E48").Copy()
$ws.Range("B49:E49").PasteSpecial(-4122)
$ws.Range("B49").Value = 45214
$ws.Range("C49").Value = 45227
$ws.Range("D49").Value = "Harian Driver"
$ws.Range("E49").Value = 21382400

$lo.ListRows.Add() | Out-Null
$ws.Range("B49:E49").Copy()
$ws.Range("B50:E50").PasteSpecial(-4122)
$ws.Range("B50").Value = 45228
$ws.Range("C50").Value = 45241
$ws.Range("D50").Value = "Harian Driver"
$ws.Range("E50").Value = 20081200

$lo.ListRows.Add() | Out-Null
$ws.Range("B50:E50").Copy()
$ws.Range("B51:E51").PasteSpecial(-4122)
$ws.Range("B51").Value = 45116
$ws.Range("C51").Value = 45129
$ws.Range("D51").Value = "Harian Fitter"
$ws.Range("E51").Value = 86520000

$lo.ListRows.Add() | Out-Null
$ws.Range("B51:E51").Copy()
$ws.Range("B52:E52").PasteSpecial(-4122)
$ws.Range("B52").Value = 45130
$ws.Range("C52").Value = 45143
$ws.Range("D52").Value = "Harian Fitter"
$ws.Range("E52").Value = 73997450

$lo.ListRows.Add() | Out-Null
$ws.Range("B52:E52").Copy()
$ws.Range("B53:E53").PasteSpecial(-4122)
$ws.Range("B53").Value = 45144
$ws.Range("C53").Value = 45157
$ws.Range("D53").Value = "Harian Fitter"
$ws.Range("E53").Value = 58252125

$lo.ListRows.Add() | Out-Null
$ws.Range("B53:E53").Copy()
$ws.Range("B54:E54").PasteSpecial(-4122)
$ws.Range("B54").Value = 45158
$ws.Range("C54").Value = 45171
$ws.Range("D54").Value = "Harian Fitter"
$ws.Range("E54").Value = 46882025

$lo.ListRows.Add() | Out-Null
$ws.Range("B54:E54").Copy()
$ws.Range("B55:E55").PasteSpecial(-4122)
$ws.Range("B55").Value = 45172
$ws.Range("C55").Value = 45185
$ws.Range("D55").Value = "Harian Fitter"
$ws.Range("E55").Value = 49916700

$lo.ListRows.Add() | Out-Null
$ws.Range("B55:E55").Copy()
$ws.Range("B56:E56").PasteSpecial(-4122)
$ws.Range("B56").Value = 45186
$ws.Range("C56").Value = 45199
$ws.Range("D56").Value = "Harian Fitter"
$ws.Range("E56").Value = 49696150

$lo.ListRows.Add() | Out-Null
$ws.Range("B56:E56").Copy()
$ws.Range("B57:E57").PasteSpecial(-4122)
$ws.Range("B57").Value = 45200
$ws.Range("C57").Value = 45213
$ws.Range("D57").Value = "Harian Fitter"
$ws.Range("E57").Value = 58555500

$lo.ListRows.Add() | Out-Null
$ws.Range("B57:E57").Copy()
$ws.Range("B58:E58").PasteSpecial(-4122)
$ws.Range("B58").Value = 45214
$ws.Range("C58").Value = 45227
$ws.Range("D58").Value = "Harian Fitter"
$ws.Range("E58").Value = 58796775

$lo.ListRows.Add() | Out-Null
$ws.Range("B58:E58").Copy()
$ws.Range("B59:E59").PasteSpecial(-4122)
$ws.Range("B59").Value = 45228
$ws.Range("C59").Value = 45241
$ws.Range("D59").Value = "Harian Fitter"
$ws.Range("E59").Value = 44149525

$lo.ListRows.Add() | Out-Null
$ws.Range("B59:E59").Copy()
$ws.Range("B60:E60").PasteSpecial(-4122)
$ws.Range("B60").Value = 45116
$ws.Range("C60").Value = 45129
$ws.Range("D60").Value = "Harian Semi Fitter, Rigger, Scaffolder Lokal"
$ws.Range("E60").Value = 54569450

$lo.ListRows.Add() | Out-Null
$ws.Range("B60:E60").Copy()
$ws.Range("B61:E61").PasteSpecial(-4122)
$ws.Range("B61").Value = 45130
$ws.Range("C61").Value = 45143
$ws.Range("D61").Value = "Harian Semi Fitter, Rigger, Scaffolder Lokal"
$ws.Range("E61").Value = 77646975

$lo.ListRows.Add() | Out-Null
$ws.Range("B61:E61").Copy()
$ws.Range("B62:E62").PasteSpecial(-4122)
$ws.Range("B62").Value = 45144
$ws.Range("C62").Value = 45157
$ws.Range("D62").Value = "Harian Semi Fitter, Rigger, Scaffolder Lokal"
$ws.Range("E62").Value = 88653300

$lo.ListRows.Add() | Out-Null
$ws.Range("B62:E62").Copy()
$ws.Range("B63:E63").PasteSpecial(-4122)
$ws.Range("B63").Value = 45158
$ws.Range("C63").Value = 45171
$ws.Range("D63").Value = "Harian Semi Fitter, Rigger, Scaffolder Lokal"
$ws.Range("E63").Value = 85171075

$lo.ListRows.Add() | Out-Null
$ws.Range("B63:E63").Copy()
$ws.Range("B64:E64").PasteSpecial(-4122)
$ws.Range("B64").Value = 45172
$ws.Range("C64").Value = 45185
$ws.Range("D64").Value = "Harian Semi Fitter, Rigger, Scaffolder Lokal"
$ws.Range("E64").Value = 80288925

$lo.ListRows.Add() | Out-Null
$ws.Range("B64:E64").Copy()
$ws.Range("B65:E65").PasteSpecial(-4122)
$ws.Range("B65").Value = 45186
$ws.Range("C65").Value = 45199
$ws.Range("D65").Value = "Harian Semi Fitter, Rigger, Scaffolder Lokal"
$ws.Range("E65").Value = 86789650

$lo.ListRows.Add() | Out-Null
$ws.Range("B65:E65").Copy()
$ws.Range("B66:E66").PasteSpecial(-4122)
$ws.Range("B66").Value = 45200
$ws.Range("C66").Value = 45213
$ws.Range("D66").Value = "Harian Semi Fitter, Rigger, Scaffolder Lokal"
$ws.Range("E66").Value = 80944575

$lo.ListRows.Add() | Out-Null
$ws.Range("B66:E66").Copy()
$ws.Range("B67:E67").PasteSpecial(-4122)
$ws.Range("B67").Value = 45214
$ws.Range("C67").Value = 45227
$ws.Range("D67").Value = "Harian Semi Fitter, Rigger, Scaffolder Lokal"
$ws.Range("E67").Value = 78798550

$lo.ListRows.Add() | Out-Null
$ws.Range("B67:E67").Copy()
$ws.Range("B68:E68").PasteSpecial(-4122)
$ws.Range("B68").Value = 45228
$ws.Range("C68").Value = 45241
$ws.Range("D68").Value = "Harian Semi Fitter, Rigger, Scaffolder Lokal"
$ws.Range("E68").Value = 69366350

$lo.ListRows.Add() | Out-Null
$ws.Range("B68:E68").Copy()
$ws.Range("B69:E69").PasteSpecial(-4122)
$ws.Range("B69").Value = 45116
$ws.Range("C69").Value = 45129
$ws.Range("D69").Value = "Harial Helper Lokal"
$ws.Range("E69").Value = 49793625

$lo.ListRows.Add() | Out-Null
$ws.Range("B69:E69").Copy()
$ws.Range("B70:E70").PasteSpecial(-4122)
$ws.Range("B70").Value = 45130
$ws.Range("C70").Value = 45143
$ws.Range("D70").Value = "Harial Helper Lokal"
$ws.Range("E70").Value = 82821750

$lo.ListRows.Add() | Out-Null
$ws.Range("B70:E70").Copy()
$ws.Range("B71:E71").PasteSpecial(-4122)
$ws.Range("B71").Value = 45144
$ws.Range("C71").Value = 45157
$ws.Range("D71").Value = "Harial Helper Lokal"
$ws.Range("E71").Value = 75367125

$lo.ListRows.Add() | Out-Null
$ws.Range("B71:E71").Copy()
$ws.Range("B72:E72").PasteSpecial(-4122)
$ws.Range("B72").Value = 45158
$ws.Range("C72").Value = 45171
$ws.Range("D72").Value = "Harial Helper Lokal"
$ws.Range("E72").Value = 87264000

$lo.ListRows.Add() | Out-Null
$ws.Range("B72:E72").Copy()
$ws.Range("B73:E73").PasteSpecial(-4122)
$ws.Range("B73").Value = 45172
$ws.Range("C73").Value = 45185
$ws.Range("D73").Value = "Harial Helper Lokal"
$ws.Range("E73").Value = 78678600

$lo.ListRows.Add() | Out-Null
$ws.Range("B73:E73").Copy()
$ws.Range("B74:E74").PasteSpecial(-4122)
$ws.Range("B74").Value = 45186
$ws.Range("C74").Value = 45199
$ws.Range("D74").Value = "Harial Helper Lokal"
$ws.Range("E74").Value = 101533500

$lo.ListRows.Add() | Out-Null
$ws.Range("B74:E74").Copy()
$ws.Range("B75:E75").PasteSpecial(-4122)
$ws.Range("B75").Value = 45200
$ws.Range("C75").Value = 45213
$ws.Range("D75").Value = "Harial Helper Lokal"
$ws.Range("E75").Value = 92418750

$lo.ListRows.Add() | Out-Null
$ws.Range("B75:E75").Copy()
$ws.Range("B76:E76").PasteSpecial(-4122)
$ws.Range("B76").Value = 45214
$ws.Range("C76").Value = 45227
$ws.Range("D76").Value = "Harial Helper Lokal"
$ws.Range("E76").Value = 91655250

$lo.ListRows.Add() | Out-Null
$ws.Range("B76:E76").Copy()
$ws.Range("B77:E77").PasteSpecial(-4122)
$ws.Range("B77").Value = 45228
$ws.Range("C77").Value = 45241
$ws.Range("D77").Value = "Harial Helper Lokal"
$ws.Range("E77").Value = 85588500

$lo.ListRows.Add() | Out-Null
$ws.Range("B77:E77").Copy()
$ws.Range("B78:E78").PasteSpecial(-4122)
$ws.Range("B78").Value = 45116
$ws.Range("C78").Value = 45129
$ws.Range("D78").Value = "Harian Mechanical Fitter"
$ws.Range("E78").Value = 4377550

$lo.ListRows.Add() | Out-Null
$ws.Range("B78:E78").Copy()
$ws.Range("B79:E79").PasteSpecial(-4122)
$ws.Range("B79").Value = 45116
$ws.Range("C79").Value = 45129
$ws.Range("D79").Value = "Harian Semi Fitter - Helper Lokal"
$ws.Range("E79").Value = 52166625

# The wider "Deskripsi" text ("Harian Semi Fitter, Rigger, Scaffolder Lokal", etc.)
# and the bigger Nominal values no longer fit the old column widths, so re-fit them.
$ws.Columns.Item(4).ColumnWidth = 35.1
$ws.Columns.Item(5).ColumnWidth = 11.6

# Author finished the session with the data sheet active and cell M72 selected
$ws2.Select()
$ws2.Range("E12").Select()
$ws.Activate()
$ws.Range("M72").Select()